# Weekly price-sheet update: a new "Arveja Verde" (Perfection, Provincia de
# Huasco) record for 2022-06-14 is inserted at the top of the data block
# (row 88), pushing every existing row below it down by one. The sheet's
# used range therefore grows from A1:R117 to A1:R118 and the last existing
# record (old row 117) ends up unchanged at row 118.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 88..117 down to 89..118, leaving row 88 free (and copying its
# previous formatting, e.g. the date style on column D, down with it).
$ws.Rows(88).Insert()

# Populate the newly freed row 88 with the new observation.
$ws.Range("A88").Value = 4
$ws.Range("B88").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C88").Value = "Los Lagos"
$ws.Range("D88").Value = 44726
$ws.Range("E88").Value = 10
$ws.Range("F88").Value = 100112022
$ws.Range("G88").Value = "Arveja Verde"
$ws.Range("H88").Value = "Perfection"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 40
$ws.Range("K88").Value = 44000
$ws.Range("L88").Value = 44000
$ws.Range("M88").Value = 44000
$ws.Range("N88").Value = "$/malla 25 kilos"
$ws.Range("O88").Value = "Provincia de Huasco"
$ws.Range("P88").Value = 1760
$ws.Range("Q88").Value = 25
$ws.Range("R88").Value = "Hortaliza"
